$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "routes" (sheet2 / Table4): add a new row describing the
# "Get Locale" route.
# ------------------------------------------------------------------
$wsRoutes = $wb.Worksheets.Item("routes")
$routesTable = $wsRoutes.ListObjects.Item("Table4")
$routesTable.ListRows.Add() | Out-Null

$wsRoutes.Range("A8:G8").WrapText = $true
$wsRoutes.Range("A8").Value = "locale"
$wsRoutes.Range("B8").Value = "Get Locale"
$wsRoutes.Range("C8").Value = "getLocale"
$wsRoutes.Range("D8").Value = "locale"
$wsRoutes.Range("F8").Value = "Gets a locale file. Returns 'en' if none specified"
$wsRoutes.Range("G8").Value = "GET"

# User was working on this sheet/cell before switching tabs.
$wsRoutes.Activate() | Out-Null
$wsRoutes.Range("F9").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "dotenvs" (sheet1 / Table1): add a new row describing the
# "default_locale" dotenv variable.
# ------------------------------------------------------------------
$wsDotenvs = $wb.Worksheets.Item("dotenvs")
$dotenvsTable = $wsDotenvs.ListObjects.Item("Table1")
$dotenvsTable.ListRows.Add() | Out-Null

$wsDotenvs.Range("A9:B9").WrapText = $true
$wsDotenvs.Range("C9").WrapText = $true
$wsDotenvs.Range("A9").Value = "Locale"
$wsDotenvs.Range("B9").Value = "default_locale"
$wsDotenvs.Range("C9").Value = "The default locale for kutenq!"
$wsDotenvs.Range("D9").Value = "en"

# dotenvs is the tab left active/selected after the edit.
$wsDotenvs.Activate() | Out-Null
$wsDotenvs.Range("D10").Select() | Out-Null
